$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("plate")

for ($i = 1; $i -le 21; $i++) {
    $row = $i + 1
    $ws.Range("A$row").Value = "20191204a$i.itc"
}
